# Updates cryptos list values (price + 1h volume change) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is unambiguously non-numeric (contains "%", multiple
# "." separators, subscript digits, etc.) can be written directly; Excel will
# keep them as text just like the original inline strings.
$ws.Range("D2").Value = '26.958.37'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '1.674.59'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("E6").Value = '  +1.58%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '1.910.38'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '1.704.10'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").Value = '26.959.86'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("E18").Value = '  +4.03%  '
$ws.Range("E19").Value = '  -1.58%  '
$ws.Range("D20").Value = '0.0₃0733'
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("E22").Value = '  -0.60%  '
$ws.Range("E23").Value = '  -1.52%  '
$ws.Range("E24").Value = '  -2.13%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("D33").Value = '1.477.87'
$ws.Range("E33").Value = '  -4.67%  '
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("E40").Value = '  -3.83%  '
$ws.Range("E41").Value = '  +5.69%  '
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("E43").Value = '  +1.85%  '
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").Value = '1.816.09'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("E51").Value = '  +0.05%  '

# Cells whose new text *looks* like a plain number (e.g. "215.01") must be
# forced to Text format first, otherwise Excel auto-converts them to a
# floating point number (losing the exact decimal text and the original
# "no explicit style" cell formatting). We flip the format to Text, assign
# the literal string, then restore the default "Normal" style so the saved
# cell matches the unstyled inline-string cell from the source workbook.
$numericTextCells = @("D5", "D11", "D16", "D19", "D22", "D23", "D25", "D28", "D35", "D37", "D40", "D44", "D46", "D47", "D51")
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D5").Value = '215.01'
$ws.Range("D11").Value = '0.0886'
$ws.Range("D16").Value = '65.64'
$ws.Range("D19").Value = '235.62'
$ws.Range("D22").Value = '4.44'
$ws.Range("D23").Value = '9.18'
$ws.Range("D25").Value = '145.48'
$ws.Range("D28").Value = '0.113'
$ws.Range("D35").Value = '1.68'
$ws.Range("D37").Value = '0.586'
$ws.Range("D40").Value = '5.86'
$ws.Range("D44").Value = '67.26'
$ws.Range("D46").Value = '0.776'
$ws.Range("D47").Value = '90.63'
$ws.Range("D51").Value = '7.73'

foreach ($ref in $numericTextCells) {
    $ws.Range($ref).Style = "Normal"
}
